# Renames the "Order" domain terminology used in the request/command
# messages of the diagram to "Address" (e.g. OrderBook -> AddressBook).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ":Order" -> ":Address"  (shape "Rectangle 62", first occurrence)
$sh = $s.Shapes.Item(6)
$tr = $sh.TextFrame.TextRange
$tr.Characters(2, 5).Text = "Address"

# "undoOrderBook()" -> "undoAddressBook()"  (shape "TextBox 78")
$sh = $s.Shapes.Item(19)
$tr = $sh.TextFrame.TextRange
$tr.Characters(5, 5).Text = "Address"

# ":VersionedOrderBook" -> ":VersionedAddressBook"  (shape "Rectangle 62")
$sh = $s.Shapes.Item(23)
$tr = $sh.TextFrame.TextRange
$tr.Characters(11, 5).Text = "Address"

# "resetData(ReadOnlyOrderBook)" -> "resetData(ReadOnlyAddressBook)"  (shape "TextBox 87")
$sh = $s.Shapes.Item(35)
$tr = $sh.TextFrame.TextRange
$tr.Characters(19, 5).Text = "Address"
